# Fruta / hortaliza, semanal
# Re-sort the weekly price records (rows 2-5) by ascending date ("Fecha").
# Columns D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) move together as a
# single record for each row; the other columns (A,B,C,E..L,Q,R,T) are
# identical for every row and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44971
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("S2").Value = 3500

$ws.Range("D3").Value = 44973
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 28000
$ws.Range("O3").Value = 28000
$ws.Range("P3").Value = 28000
$ws.Range("S3").Value = 3500

$ws.Range("D4").Value = 44981
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 25000
$ws.Range("S4").Value = 3125

$ws.Range("D5").Value = 44980
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 25000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 25000
$ws.Range("S5").Value = 3125
